# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Angeleno" (Ciruela) at row 207,
# shifting the existing rows 207:301 down to 209:303.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 207 (this pushes the old
# rows 207-301 down to 209-303 and extends the used range to A1:T303).
$ws.Rows("207:208").Insert()

# --- Row 207: Angeleno / Primera ---
$ws.Range("A207").Value = 4
$ws.Range("B207").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C207").Value = "Los Lagos"
$ws.Range("D207").Value = 44992
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100103
$ws.Range("H207").Value = "Frutos de hueso (carozo)"
$ws.Range("I207").Value = 100103002
$ws.Range("J207").Value = "Ciruela"
$ws.Range("K207").Value = "Angeleno"
$ws.Range("L207").Value = "Primera"
$ws.Range("M207").Value = 400
$ws.Range("N207").Value = 15000
$ws.Range("O207").Value = 16000
$ws.Range("P207").Value = 15500
$ws.Range("Q207").Value = "$/caja 14 kilos granel"
$ws.Range("R207").Value = "Región Metropolitana"
$ws.Range("S207").Value = 1107
$ws.Range("T207").Value = 14

# --- Row 208: Angeleno / Segunda ---
$ws.Range("A208").Value = 4
$ws.Range("B208").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C208").Value = "Los Lagos"
$ws.Range("D208").Value = 44992
$ws.Range("E208").Value = 10
$ws.Range("F208").Value = "Fruta"
$ws.Range("G208").Value = 100103
$ws.Range("H208").Value = "Frutos de hueso (carozo)"
$ws.Range("I208").Value = 100103002
$ws.Range("J208").Value = "Ciruela"
$ws.Range("K208").Value = "Angeleno"
$ws.Range("L208").Value = "Segunda"
$ws.Range("M208").Value = 200
$ws.Range("N208").Value = 14000
$ws.Range("O208").Value = 14000
$ws.Range("P208").Value = 14000
$ws.Range("Q208").Value = "$/caja 14 kilos granel"
$ws.Range("R208").Value = "Región Metropolitana"
$ws.Range("S208").Value = 1000
$ws.Range("T208").Value = 14
